$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.566.63"
$ws.Range("E2").Value = "  +1.95%  "
$ws.Range("D3").Value = "1.580.21"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").Value = "'212.62"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "'0.492"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").Value = "'46.95"
$ws.Range("E8").Value = "  +8.20%  "
$ws.Range("D9").Value = "'24.12"
$ws.Range("E9").Value = "  +3.87%  "
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").Value = "'0.0592"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").Value = "'0.0882"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "1.804.54"
$ws.Range("D14").Value = "1.587.74"
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("D15").Value = "'0.525"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").Value = "'3.71"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").Value = "28.562.42"
$ws.Range("E17").Value = "  +1.94%  "
$ws.Range("D18").Value = "'62.44"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D19").Value = "'229.25"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").Value = "'7.43"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("E21").Value = "  -1.37%  "
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("E23").Value = "  -4.48%  "
$ws.Range("D24").Value = "'9.17"
$ws.Range("E24").Value = "  -1.56%  "
$ws.Range("E25").Value = "  +5.68%  "
$ws.Range("D26").Value = "'151.16"
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("D27").Value = "'15.04"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("E28").Value = "  -1.51%  "
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("E31").Value = "  -1.55%  "
$ws.Range("D32").Value = "'0.0464"
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("D33").Value = "'3.22"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "1.400.94"
$ws.Range("E35").Value = "  -1.14%  "
$ws.Range("D36").Value = "'1.56"
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("E37").Value = "  -2.14%  "
$ws.Range("E38").Value = "  +2.04%  "
$ws.Range("E39").Value = "  +6.20%  "
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").Value = "'0.531"
$ws.Range("E41").Value = "  -1.49%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'0.798"
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.51%  "
$ws.Range("D44").Value = "'5.62"
$ws.Range("E44").Value = "  -0.75%  "
$ws.Range("D45").Value = "'1.86"
$ws.Range("E45").Value = "  +2.20%  "
$ws.Range("E46").Value = "  +0.67%  "
$ws.Range("D47").Value = "'62.77"
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("D48").Value = "1.715.63"
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("D49").Value = "'86.03"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("E50").Value = "  -2.37%  "
$ws.Range("D51").Value = "'0.0517"
$ws.Range("E51").Value = "  -1.16%  "
